$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "RelivePos" column (G) previously held distinct unique placeholder
# values (0,0,0 .. 0,0,6) for the revive/move message body. The unused
# per-row variants are removed and every row now shares the single
# updated body value.
$ws.Range("G11:G17").Value = "55,110,0"

# Update the active selection to match the authored state (single cell).
$ws.Range("G15").Select()
